# Doyle + LG_M50 Buildable - Cleanup
# Insert a new title row at the very top of the sheet (everything else shifts
# down by one row, formulas auto-adjust), label it, then tidy up by clearing
# a stray leftover formula and stamping out a block of blank, but still
# formatted, rows further down the sheet (mirrors the original formatted
# block) ready for the next data set to be pasted in underneath.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing table down one row and add the new heading.
$ws.Rows("1:1").Insert()
$ws.Range("A1").Value = "Chen 2020"

# The old stray "0.000173" helper formula (now at F11) is no longer needed,
# clear its contents but keep its number formatting.
$ws.Range("F11").ClearContents()

# Re-stamp the formatting of the original calculation block 19 rows further
# down the sheet, ready for another data set - values/formulas are left
# blank, only the cell formatting is carried across.
$ws.Range("A3:D3").Copy()
$ws.Range("A22:D22").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("F22").PasteSpecial(-4122)

$ws.Range("A6").Copy()
$ws.Range("A25").PasteSpecial(-4122)

$ws.Range("B10:F10").Copy()
$ws.Range("B29:F29").PasteSpecial(-4122)

$ws.Range("F11").Copy()
$ws.Range("F30").PasteSpecial(-4122)

$ws.Range("A16").Copy()
$ws.Range("A35").PasteSpecial(-4122)

# Clear the clipboard marching ants / final selection used for editing.
$ws.Range("A2").Select()
